$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = ""
$ws.Range("B31").Value = "يامن "
$ws.Range("C31").Value = "22"
$ws.Range("D31").Value = "الصمود"
$ws.Range("E31").Value = "الرحلة 1"
$ws.Range("F31").Value = "C3"
$ws.Range("G31").Value = "NRC"
$ws.Range("H31").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٠٩:٥٦ م"
